# feat: add 2022-Q3 data
#
# Before: sheet1 = "总计" (summary), sheet2 = "2022-Q2" (fund holdings for Q2)
# After:  sheet1 = "总计" (summary, +1 row for Q3),
#         sheet2 = "2022-Q3" (NEW fund holdings data, reuses the old "2022-Q2" sheet/rId),
#         sheet3 = "2022-Q2" (the old fund holdings data, moved to a brand-new sheet,
#                              keeping its original look)

$wb = $excel.ActiveWorkbook
$wsTotal = $wb.Worksheets.Item(1)
$wsQ2 = $wb.Worksheets.Item(2)

# A never-touched cell used purely as a "default style" donor for
# PasteSpecial(formats) resets (keeps text values while clearing any
# quote-prefix style picked up from forcing text with a leading apostrophe).
$blankCell = $wsTotal.Range("Z100")

# Stash the old "2022-Q2" sheet's header-row / column-A look (style "1" in
# the source workbook - distinct from the "总计" sheet's style "2") onto
# holding cells before anything else changes, so it can be reapplied to the
# relocated sheet further down even after the source sheet is overwritten.
$wsQ2.Range("B1").Copy()
$wsTotal.Range("Z90").PasteSpecial(-4122)
$wsQ2.Range("A2").Copy()
$wsTotal.Range("Z91").PasteSpecial(-4122)
$oldQ2HeaderStyle = $wsTotal.Range("Z90")
$oldQ2AColStyle = $wsTotal.Range("Z91")

# ---------------------------------------------------------------------------
# 1) Capture the existing "2022-Q2" fund-holdings rows before we overwrite
#    them - they get relocated verbatim onto a brand-new worksheet.
# ---------------------------------------------------------------------------
$oldQ2Rows = @(
    @("002558", "博时鑫瑞灵活配置混合A", "4.31", "22.12", "1.12", "0.0483", 5),
    @("003300", "华夏圆和灵活配置混合", "0.58", "81.65", "5.27", "0.0306", 4),
    @("002559", "博时鑫瑞灵活配置混合C", "1.26", "22.12", "1.12", "0.0141", 5),
    @("015068", "华夏圆和灵活配置混合C", "0.03", "81.65", "5.27", "0.0016", 4)
)

# ---------------------------------------------------------------------------
# 2) "总计" sheet: insert a new row 2 for 2022-Q3, push the 2022-Q2 row to
#    row 3 (and bump its running index from 0 to 1).
# ---------------------------------------------------------------------------
$oldB2 = $wsTotal.Cells.Item(2, 2).Value()
$oldC2 = $wsTotal.Cells.Item(2, 3).Value()
$oldD2 = $wsTotal.Cells.Item(2, 4).Value()

$wsTotal.Cells.Item(3, 1).Value = 1
$wsTotal.Cells.Item(3, 2).Value = $oldB2
$wsTotal.Cells.Item(3, 3).Value = $oldC2
$wsTotal.Cells.Item(3, 4).Value = $oldD2

$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q3"
$wsTotal.Cells.Item(2, 3).Value = 8
$wsTotal.Cells.Item(2, 4).Value = 0.87

# Re-apply the A-column style (bold/border/center) to the row-3 cell so it
# matches row 2 and the sheet's existing look.
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) The sheet that is currently "2022-Q2" becomes "2022-Q3": rename it and
#    replace its contents with the Q3 fund-holdings table.
# ---------------------------------------------------------------------------
$wsQ2.Cells.Clear()
$wsQ3 = $wsQ2
$wsQ3.Name = "2022-Q3"

$q3Header = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$col = 2
foreach ($h in $q3Header) {
    $wsQ3.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}

$q3Rows = @(
    @("004640", "华夏节能环保股票A", "6.04", "93.11", "4.68", "0.2827", 7),
    @("012703", "华夏核心成长混合A", "4.94", "93.34", "4.88", "0.2411", 6),
    @("015229", "华夏低碳经济一年持有混合A", "2.49", "91.97", "6.57", "0.1636", 4),
    @("015230", "华夏低碳经济一年持有混合C", "1.02", "91.97", "6.57", "0.0670", 4),
    @("003300", "华夏圆和灵活配置混合A", "0.77", "75.31", "6.99", "0.0538", 3),
    @("012710", "华夏核心成长混合C", "0.61", "93.34", "4.88", "0.0298", 6),
    @("015068", "华夏圆和灵活配置混合C", "0.33", "75.31", "6.99", "0.0231", 3),
    @("015060", "华夏节能环保股票C", "0.26", "93.11", "4.68", "0.0122", 7)
)

$r = 2
foreach ($row in $q3Rows) {
    $wsQ3.Cells.Item($r, 1).Value = $r - 2
    $wsQ3.Cells.Item($r, 2).Value = "'" + $row[0]
    $wsQ3.Cells.Item($r, 3).Value = $row[1]
    $wsQ3.Cells.Item($r, 4).Value = "'" + $row[2]
    $wsQ3.Cells.Item($r, 5).Value = "'" + $row[3]
    $wsQ3.Cells.Item($r, 6).Value = "'" + $row[4]
    $wsQ3.Cells.Item($r, 7).Value = "'" + $row[5]
    $wsQ3.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# Clear the stray "quote prefix" style the leading apostrophes above created
# (text values are kept - PasteSpecial(formats) never touches cell content).
$blankCell.Copy()
$wsQ3.Range("B2:G9").PasteSpecial(-4122)

# Style the new sheet like the "总计" sheet: header row + column A use the
# bold/border/center style ("s=2" in the source workbook).
$wsTotal.Range("B1").Copy()
$wsQ3.Range("B1:H1").PasteSpecial(-4122)
$wsTotal.Range("A2").Copy()
$wsQ3.Range("A2:A9").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4) Add a brand-new worksheet right after "2022-Q3" to hold the relocated
#    "2022-Q2" fund-holdings data (matching the old sheet's original look).
# ---------------------------------------------------------------------------
$wsQ2New = $wb.Worksheets.Add($null, $wsQ3)
$wsQ2New.Name = "2022-Q2"

$col = 2
foreach ($h in $q3Header) {
    $wsQ2New.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}

$r = 2
foreach ($row in $oldQ2Rows) {
    $wsQ2New.Cells.Item($r, 1).Value = $r - 2
    $wsQ2New.Cells.Item($r, 2).Value = "'" + $row[0]
    $wsQ2New.Cells.Item($r, 3).Value = $row[1]
    $wsQ2New.Cells.Item($r, 4).Value = "'" + $row[2]
    $wsQ2New.Cells.Item($r, 5).Value = "'" + $row[3]
    $wsQ2New.Cells.Item($r, 6).Value = "'" + $row[4]
    $wsQ2New.Cells.Item($r, 7).Value = "'" + $row[5]
    $wsQ2New.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# Clear the stray "quote prefix" style from the leading apostrophes above.
$blankCell.Copy()
$wsQ2New.Range("B2:G5").PasteSpecial(-4122)

# Restore the old "2022-Q2" sheet's original header-row / column-A style
# (style "1" in the source workbook) from the holding cells stashed earlier.
$oldQ2HeaderStyle.Copy()
$wsQ2New.Range("B1:H1").PasteSpecial(-4122)
$oldQ2AColStyle.Copy()
$wsQ2New.Range("A2:A5").PasteSpecial(-4122)

# Clean up the holding cells so they don't leak into "总计"'s used range.
$wsTotal.Range("Z90:Z91").Clear()

$wsTotal.Select()
$wsTotal.Range("A1").Select()

Write-Output "done"
